$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the existing drug name to lowercase ("Crocin" -> "crocin")
$ws.Range("A2").Value = "crocin"

# Add a new "Params" column to the database schema
$ws.Range("C1").Value = "Params"
$ws.Range("C2").Value = "ex"

# Touch D1 so it participates in the sheet's used range/style like the rest
# of the header row (no content, just formatting/presence)
$ws.Range("D1").Style = "Normal"

# Insert a new record row
$ws.Range("A3").Value = "crocin"
$ws.Range("B3").Value = 0
$ws.Range("D3").Value = "g"

# Leave the selection where the editing session ended up
$ws.Range("F4").Select() | Out-Null
